$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("Users will be displayed their new details.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rng.Find.Found) {
    throw "Anchor text not found"
}

# Find which paragraph holds the matched text, so we can insert a brand-new
# paragraph right after it.
$paras = $d.Paragraphs
$anchorIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*Users will be displayed their new details.*") {
        $anchorIndex = $i
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not locate anchor paragraph"
}

$rng.Collapse(0)
$rng.InsertParagraphAfter()

$paras = $d.Paragraphs
$newPara = $paras.Item($anchorIndex + 1)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Create Attendance Sheets</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Times-Bold"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Description: </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>Admins can create new attendance sheets</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Times-Bold"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Actor: </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>User(Admins)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Times-Bold"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Input: </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Admins provide list of </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>UIDs ,day</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> range and mode of attendance.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Times-Bold"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>Output</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:cs="Times-Bold"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>Admins</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> will be given an </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>acknoledgemend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cs="Times-Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and taken to their attendance sheet.</w:t></w:r></w:p>'
$newPara.Range.InsertXML($xml)
Write-Host "Insertion complete. Paragraph count: $($d.Paragraphs.Count)"
